# feat: add 2022-Q4 data
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q4" worksheet before the existing "2022-Q2"
#    sheet, so the tab order becomes: 总计, 2022-Q4, 2022-Q2, 2021-Q4
# ------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q4Sheet = $wb.Worksheets.Add($q2Sheet)
$q4Sheet.Name = "2022-Q4"

# Header row for the new sheet (copy formatting from "2022-Q2"'s header,
# then re-apply the actual header text to every cell in that range so the
# pasted formatting is committed for each individual cell)
$q2Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# Data row for the new sheet (copy formatting for the index cell A2)
$q2Sheet.Range("A2").Copy()
$q4Sheet.Range("A2").PasteSpecial(-4122)
$q4Sheet.Range("A2").Value = 0

# Fund code / fund name stay as plain text
$q4Sheet.Range("B2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "159620"
$q4Sheet.Range("C2").Value = "华夏中证智选500成长创新策略ETF"

# Numeric-looking figures are stored as text (matches the other quarter
# sheets, which keep the original fixed-decimal formatting, e.g. "0.0050")
$q4Sheet.Range("D2").NumberFormat = "@"
$q4Sheet.Range("D2").Value = "0.32"
$q4Sheet.Range("E2").NumberFormat = "@"
$q4Sheet.Range("E2").Value = "95.01"
$q4Sheet.Range("F2").NumberFormat = "@"
$q4Sheet.Range("F2").Value = "1.57"
$q4Sheet.Range("G2").NumberFormat = "@"
$q4Sheet.Range("G2").Value = "0.0050"
$q4Sheet.Range("H2").Value = 8

# ------------------------------------------------------------------
# 2. Add a summary row for "2022-Q4" on the "总计" sheet, right
#    after the header row (so it becomes the new row 2, pushing the
#    previous rows 2 & 3 down to rows 3 & 4).
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").Style = "Normal"

# Restore the index-column (A) formatting to match the other rows, then
# write every cell of the new row
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0

# Renumber the index column so it keeps reading 0,1,2
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
